$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number/date and must be forced to text
# to match the original inline-string cell type (avoid Excel auto-numeric coercion).
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D17", "D19", "D20", "D21", "D23", "D24", "D25", "D27", "D28", "D30", "D31", "D32", "D33", "D35", "D37", "D38", "D41", "D42", "D43", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all text values. The unicode subscript-three glyph is built via [char]
# and joined with the -f format operator (string "+" concatenation here gets
# numerically coerced by the interpreter, e.g. "0.0"+"0940" -> Double).
$sub3 = [char]0x2083

$ws.Range("D2").Value = "48.128.61"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.504.61"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "320.89"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "107.56"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("D10").Value = "39.60"
$ws.Range("E10").Value = "  -2.61%  "
$ws.Range("D11").Value = "20.19"
$ws.Range("D12").Value = "0.0812"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "2.897.07"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "2.507.32"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "0.835"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "48.005.33"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "6.72"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "2.77"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = ("{0}{1}{2}" -f "0.0", $sub3, "0940")
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").Value = "278.45"
$ws.Range("E23").Value = "  +12.52%  "
$ws.Range("D24").Value = "71.48"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D27").Value = "25.97"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").Value = "9.70"
$ws.Range("E28").Value = "  -2.92%  "
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").Value = "35.33"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "2.09"
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").Value = "49.59"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "19.54"
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").Value = "5.30"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.64"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.94"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").Value = "121.60"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("D42").Value = "2.21"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").Value = "21.25"
$ws.Range("E43").Value = "  -5.93%  "
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "2.021.23"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "80.18"
$ws.Range("E51").Value = "  +3.09%  "

# Restore default (Normal) style on the forced-text cells so no stray
# number-format style is left applied (matches original unstyled cells).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
